# Add Arabic ("ara") translations of the id_type master-data rows.
# Mirrors the existing "eng" (rows 2-6) / "fra" (rows 7-11) blocks by
# appending 5 new rows (12-16) for lang_code "ara".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# code / name / descr for each of the 5 existing id types, translated to Arabic.
$rows = @(
    @{ Row = 12; Code = "UIN";      Name = "رقم التعريف الفريد";                  Descr = "الهوية الوطنية الممنوحة لمقدم الطلب" },
    @{ Row = 13; Code = "PRID";     Name = "معرف التسجيل المسبق";                 Descr = "تم تعيين المعرف بعد التسجيل المسبق" },
    @{ Row = 14; Code = "RID";      Name = "معرف تسجيل";                         Descr = "المعرف المعين بعد التسجيل" },
    @{ Row = 15; Code = "VID";      Name = "المعرف الظاهري";                      Descr = "المعرف المستخدم في استبدال UIN" },
    @{ Row = 16; Code = "Token ID"; Name = "معرف الرمز";                          Descr = "المعرف المستخدم من قبل البائع لمقدم الطلب" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # lang_code / code columns (A, B) - plain text, reuses existing shared strings.
    $ws.Cells.Item($rowNum, 1).Value = "ara"
    $ws.Cells.Item($rowNum, 2).Value = $r.Code

    # name / descr columns (C, D) - left aligned + wrapped, like a long text field.
    $ws.Cells.Item($rowNum, 3).Value = $r.Name
    $ws.Cells.Item($rowNum, 4).Value = $r.Descr
    $cd = $ws.Range("C" + $rowNum + ":D" + $rowNum)
    $cd.HorizontalAlignment = -4131
    $cd.WrapText = $true

    # is_active column (E) - copy the existing "TRUE" text cell so it stays a
    # text value instead of being coerced into a boolean.
    $ws.Range("E11").Copy($ws.Range("E" + $rowNum)) | Out-Null

    $ws.Rows.Item($rowNum).RowHeight = 16.4
}

# Column widths to fit the new Arabic text (approximate Excel's character-width
# rounding to land as close as possible to the authored values).
$ws.Columns.Item(1).ColumnWidth = 11.09
$ws.Columns.Item(2).ColumnWidth = 10.75
$ws.Columns.Item(3).ColumnWidth = 36.75
$ws.Columns.Item(4).ColumnWidth = 44.25

# Match the selection left behind by the edit.
$ws.Range("C13:D16").Select() | Out-Null
